$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Capture existing comments (text) before we start shifting rows around,
#    so we can re-create them at their new locations afterwards.
$commentB12Text = $ws.Range("B12").Comment.Text()
$commentB15Text = $ws.Range("B15").Comment.Text()

# Remove the original comments - they will be re-added at the shifted
# locations (B13 / B16) once the new row has been inserted below.
$ws.Range("B12").Comment.Delete()
$ws.Range("B15").Comment.Delete()

# 2. Insert a brand new row above row 3 (pushes rows 3-19 down to 4-20),
#    picking up the banding style ("s=2") from row 2 automatically.
$ws.Rows("3:3").Insert()

# 3. Populate the newly inserted row with the new task.
$ws.Range("A3").Value = "Engine"
$ws.Range("B3").Value = "Refactor index based unique container.  Think about using size_t rather than unsinged int."
$ws.Range("C3").Value = 2

# 4. "Basics of a physics engine" (now on row 9 after the insert above)
#    has had its estimate revised from 21 to 35.
$ws.Range("C9").Value = 35

# 5. The very last row (now row 20, "Scene Exporter" / "Get 64-bit scene
#    exporter working") has been removed completely.
$ws.Rows("20:20").Delete()

# 6. Re-create the two comments at their shifted locations (the row insert
#    above moved the underlying tasks down by one row, from B12->B13 and
#    B15->B16), preserving their original text.
$ws.Range("B13").AddComment($commentB12Text)
$ws.Range("B16").AddComment($commentB15Text)

# 7. Update the active cell / selection recorded in the sheet view.
$ws.Range("B18").Select()
